# cryptos.xlsx refresh -- GitHub Actions run on Sun Sep 17 04:31:51 UTC 2023.
# Column D (Price) and column E (Volume(1h)) are refreshed for every coin row;
# rows 48/49 (BabyDogeCoin / Cronos) also swapped places in the ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: cell reference, new text, and whether the new text had to be
# quote-prefixed below to stop it being auto-parsed as a number (the source
# file stores every Price/Volume cell as text, even fully numeric prices).
$updates = @(
    @("D2", '26.832.69', $false),
    @("E2", '  +0.42%  ', $false),
    @("D3", '1.646.67', $false),
    @("E3", '  -0.09%  ', $false),
    @("E4", '  +0.53%  ', $false),
    @("D5", '''217.38', $true),
    @("E5", '  +0.62%  ', $false),
    @("E6", '  -0.29%  ', $false),
    @("E7", '  +0.64%  ', $false),
    @("E8", '  -0.70%  ', $false),
    @("E9", '  -0.11%  ', $false),
    @("D10", '''19.21', $true),
    @("E10", '  -0.85%  ', $false),
    @("D11", '''0.0843', $true),
    @("E11", '  -0.14%  ', $false),
    @("D12", '1.870.51', $false),
    @("E12", '  -0.41%  ', $false),
    @("D13", '1.643.98', $false),
    @("E13", '  -0.27%  ', $false),
    @("D14", '''4.18', $true),
    @("E14", '  -1.04%  ', $false),
    @("D15", '''0.528', $true),
    @("E15", '  -1.17%  ', $false),
    @("D16", '''64.59', $true),
    @("E16", '  -2.70%  ', $false),
    @("D17", '26.818.65', $false),
    @("E17", '  +0.20%  ', $false),
    @("D18", '0.0₃0738', $false),
    @("E18", '  -2.23%  ', $false),
    @("D19", '''214.26', $true),
    @("E19", '  -2.70%  ', $false),
    @("E20", '  +0.68%  ', $false),
    @("D21", '''4.37', $true),
    @("E21", '  -0.44%  ', $false),
    @("D22", '''2.41', $true),
    @("E22", '  +13.18%  ', $false),
    @("D23", '''6.29', $true),
    @("E23", '  -0.72%  ', $false),
    @("D24", '''9.38', $true),
    @("E24", '  -1.97%  ', $false),
    @("D25", '''145.43', $true),
    @("E25", '  -1.16%  ', $false),
    @("E26", '  +0.75%  ', $false),
    @("D27", '''0.119', $true),
    @("E27", '  -1.77%  ', $false),
    @("E28", '  -0.06%  ', $false),
    @("D29", '''15.70', $true),
    @("E29", '  -1.27%  ', $false),
    @("D30", '''0.0514', $true),
    @("E30", '  -0.95%  ', $false),
    @("E31", '  +0.11%  ', $false),
    @("D32", '''3.33', $true),
    @("E32", '  -2.66%  ', $false),
    @("D33", '''3.01', $true),
    @("E33", '  -1.64%  ', $false),
    @("D34", '1.292.90', $false),
    @("E34", '  +0.41%  ', $false),
    @("E35", '  -1.10%  ', $false),
    @("E36", '  +1.44%  ', $false),
    @("D37", '''0.0175', $true),
    @("E37", '  -5.12%  ', $false),
    @("D38", '''0.538', $true),
    @("E38", '  +2.43%  ', $false),
    @("D39", '''0.828', $true),
    @("E39", '  -0.37%  ', $false),
    @("E40", '  +0.68%  ', $false),
    @("D41", '''0.810', $true),
    @("E41", '  -0.35%  ', $false),
    @("E42", '  -0.29%  ', $false),
    @("D43", '''5.37', $true),
    @("E43", '  -1.22%  ', $false),
    @("D44", '1.795.77', $false),
    @("E44", '  +0.42%  ', $false),
    @("D45", '''60.64', $true),
    @("E45", '  +1.45%  ', $false),
    @("D46", '''91.48', $true),
    @("E46", '  -2.35%  ', $false),
    @("E47", '  -0.04%  ', $false),
    @("B48", 'Cronos', $false),
    @("C48", 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', $false),
    @("D48", '''0.0521', $true),
    @("E48", '  +0.87%  ', $false),
    @("B49", 'BabyDogeCoin', $false),
    @("C49", 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge', $false),
    @("D49", '0.0₆0101', $false),
    @("E49", '  -4.41%  ', $false),
    @("D50", '''7.69', $true),
    @("E50", '  -1.28%  ', $false),
    @("D51", '''0.0980', $true),
    @("E51", '  +0.19%  ', $false)
)

foreach ($update in $updates) {
    $cellRef = $update[0]
    $newValue = $update[1]
    $quotePrefixed = $update[2]
    $ws.Range($cellRef).Value = $newValue
    if ($quotePrefixed) {
        # Drop the text-number-format side effect of the quote prefix so the
        # cell keeps its original (default) style, same as every other cell.
        $ws.Range($cellRef).Style = "Normal"
    }
}

